$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price figures that were authored as plain text (note some
# use "."-grouped thousands, e.g. "29.329.46"). For the cells whose new value
# would otherwise be auto-recognised by Excel as a genuine number, force the
# cell to Text format first so the stored value keeps its original string type.

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.329.67"
$ws.Range("D3").Value = "1.873.25"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "0.7119"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "241.86"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.3108"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07777"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("D10").Value = "25.10"
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("D11").Value = "0.08399"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "1.868.38"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "5.243"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "0.7107"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "91.13"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "29.334.86"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "6.072"
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("D18").Value = "0.000008186"
$ws.Range("E18").Value = "  +4.90%  "
$ws.Range("D19").Value = "240.00"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("D20").Value = "13.19"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "2.121.81"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "7.757"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "0.1600"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").Value = "162.72"
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").Value = "9.025"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "1.508"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").Value = "4.408"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "4.319"
$ws.Range("E31").Value = "  +1.33%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "1.287"
$ws.Range("E32").Value = "  -3.10%  "
$ws.Range("D33").Value = "0.05289"
$ws.Range("E33").Value = "  +2.87%  "
$ws.Range("D34").Value = "1.937"
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("D35").Value = "1.176"
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("D36").Value = "0.7463"
$ws.Range("E36").Value = "  -6.26%  "
$ws.Range("D37").Value = "2.702"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").Value = "0.01873"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("D39").Value = "1.220.13"
$ws.Range("E39").Value = "  +4.68%  "
$ws.Range("D40").Value = "2.724"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").Value = "6.532"
$ws.Range("E41").Value = "  +5.18%  "
$ws.Range("D42").Value = "109.86"
$ws.Range("E42").Value = "  +7.62%  "
$ws.Range("D43").Value = "0.8870"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").Value = "72.41"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "2.019.23"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "1.799"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").Value = "0.5194"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("D49").Value = "0.00000000123"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").Value = "9.373"
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").Value = "0.4312"
$ws.Range("E51").Value = "  +0.94%  "
